$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed/modified) date column C for every data row
# (rows 2-91) from 45177 (2023-09-08) to 45178 (2023-09-09).
for ($r = 2; $r -le 91; $r++) {
    $ws.Cells.Item($r, 3).Value = 45178
}

# Row 2 also lost the species "Porslinsblå spindling" (a VU / threatened /
# red-listed species), so update the related counts and the species list.
$ws.Range("K2").Value = 0   # VU
$ws.Range("O2").Value = 5   # Rödlistade (red-listed)
$ws.Range("P2").Value = 0   # Hotade (threatened)
$ws.Range("Q2").Value = 19  # Alla arter (all species)

$species = @(
    "Barrviolspindling",
    "Dvärgbägarlav",
    "Leptoporus erubescens",
    "Spillkråka",
    "Vedtrappmossa",
    "Bronshjon",
    "Dropptaggsvamp",
    "Fällmossa",
    "Grön sköldmossa",
    "Guldlockmossa",
    "Kornknutmossa",
    "Rödgul trumpetsvamp",
    "Stubbspretmossa",
    "Svavelriska",
    "Sårläka",
    "Vågbandad barkbock",
    "Vanlig groda",
    "Blåsippa",
    "Revlummer"
)
$ws.Range("R2").Value = ($species -join "`r`n")
